$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 1073-1075: header rows of sentence 51 (C=51), red-filled style (s=2) ---
$ws.Cells.Item(1073, 3).Value = 51
$ws.Cells.Item(1073, 4).Value = 0
$ws.Cells.Item(1074, 4).Value = 1
$ws.Cells.Item(1075, 4).Value = 2
$ws.Range("D1073:F1075").Interior.Color = 255

# --- Sentence 51 words: rows 1076-1092 (D=3..19), column E=word, F=count ---
$ws.Cells.Item(1076, 4).Value = 3
$ws.Cells.Item(1076, 5).Value = "[b'however']"
$ws.Cells.Item(1076, 6).Value = 1
$ws.Cells.Item(1077, 4).Value = 4
$ws.Cells.Item(1077, 5).Value = "[b'leading']"
$ws.Cells.Item(1077, 6).Value = 1
$ws.Cells.Item(1078, 4).Value = 5
$ws.Cells.Item(1078, 5).Value = "[b'questions']"
$ws.Cells.Item(1078, 6).Value = 1
$ws.Cells.Item(1079, 4).Value = 6
$ws.Cells.Item(1079, 5).Value = "[b'suggest']"
$ws.Cells.Item(1079, 6).Value = 1
$ws.Cells.Item(1080, 4).Value = 7
$ws.Cells.Item(1080, 5).Value = "[b'an']"
$ws.Cells.Item(1080, 6).Value = 1
$ws.Cells.Item(1081, 4).Value = 8
$ws.Cells.Item(1081, 5).Value = "[b'answer']"
$ws.Cells.Item(1081, 6).Value = 1
$ws.Cells.Item(1082, 4).Value = 9
$ws.Cells.Item(1082, 5).Value = "[b'and']"
$ws.Cells.Item(1082, 6).Value = 1
$ws.Cells.Item(1083, 4).Value = 10
$ws.Cells.Item(1083, 5).Value = "[b'have']"
$ws.Cells.Item(1083, 6).Value = 1
$ws.Cells.Item(1084, 4).Value = 11
$ws.Cells.Item(1084, 5).Value = "[b'the']"
$ws.Cells.Item(1084, 6).Value = 1
$ws.Cells.Item(1085, 4).Value = 12
$ws.Cells.Item(1085, 5).Value = "[b'potential']"
$ws.Cells.Item(1085, 6).Value = 1
$ws.Cells.Item(1086, 4).Value = 13
$ws.Cells.Item(1086, 5).Value = "[b'to']"
$ws.Cells.Item(1086, 6).Value = 1
$ws.Cells.Item(1087, 4).Value = 14
$ws.Cells.Item(1087, 5).Value = "[b'decrease']"
$ws.Cells.Item(1087, 6).Value = 1
$ws.Cells.Item(1088, 4).Value = 15
$ws.Cells.Item(1088, 5).Value = "[b'the']"
$ws.Cells.Item(1088, 6).Value = 1
$ws.Cells.Item(1089, 4).Value = 16
$ws.Cells.Item(1089, 5).Value = "[b'accuracy']"
$ws.Cells.Item(1089, 6).Value = 1
$ws.Cells.Item(1090, 4).Value = 17
$ws.Cells.Item(1090, 5).Value = "[b'of']"
$ws.Cells.Item(1090, 6).Value = 1
$ws.Cells.Item(1091, 4).Value = 18
$ws.Cells.Item(1091, 5).Value = "[b'the']"
$ws.Cells.Item(1091, 6).Value = 1
$ws.Cells.Item(1092, 4).Value = 19
$ws.Cells.Item(1092, 5).Value = "[b'answer']"
$ws.Cells.Item(1092, 6).Value = 1

# --- Sentence 52: row 1094 starts with C=52 (green style s=4), rows 1094-1113 (D=0..19) ---
$ws.Cells.Item(1094, 3).Value = 52
$ws.Cells.Item(1094, 4).Value = 0
$ws.Cells.Item(1094, 5).Value = "[b'made']"
$ws.Cells.Item(1094, 6).Value = 1
$ws.Cells.Item(1095, 4).Value = 1
$ws.Cells.Item(1095, 5).Value = "[b'by']"
$ws.Cells.Item(1095, 6).Value = 1
$ws.Cells.Item(1096, 4).Value = 2
$ws.Cells.Item(1096, 5).Value = "[b'new']"
$ws.Cells.Item(1096, 6).Value = 1
$ws.Cells.Item(1097, 4).Value = 3
$ws.Cells.Item(1097, 5).Value = "[b'interviewers']"
$ws.Cells.Item(1097, 6).Value = 1
$ws.Cells.Item(1098, 4).Value = 4
$ws.Cells.Item(1098, 5).Value = "[b'is']"
$ws.Cells.Item(1098, 6).Value = 1
$ws.Cells.Item(1099, 4).Value = 5
$ws.Cells.Item(1099, 5).Value = "[b'asking']"
$ws.Cells.Item(1099, 6).Value = 1
$ws.Cells.Item(1100, 4).Value = 6
$ws.Cells.Item(1100, 5).Value = "[b'too']"
$ws.Cells.Item(1100, 6).Value = 1
$ws.Cells.Item(1101, 4).Value = 7
$ws.Cells.Item(1101, 5).Value = "[b'many']"
$ws.Cells.Item(1101, 6).Value = 1
$ws.Cells.Item(1102, 4).Value = 8
$ws.Cells.Item(1102, 5).Value = "[b'narrow']"
$ws.Cells.Item(1102, 6).Value = 1
$ws.Cells.Item(1103, 4).Value = 9
$ws.Cells.Item(1103, 5).Value = "[b'questions']"
$ws.Cells.Item(1103, 6).Value = 1
$ws.Cells.Item(1104, 4).Value = 10
$ws.Cells.Item(1104, 5).Value = "[b'and']"
$ws.Cells.Item(1104, 6).Value = 1
$ws.Cells.Item(1105, 4).Value = 11
$ws.Cells.Item(1105, 5).Value = "[b'the']"
$ws.Cells.Item(1105, 6).Value = 1
$ws.Cells.Item(1106, 4).Value = 12
$ws.Cells.Item(1106, 5).Value = "[b'question']"
$ws.Cells.Item(1106, 6).Value = 1
$ws.Cells.Item(1107, 4).Value = 13
$ws.Cells.Item(1107, 5).Value = "[b'elicits']"
$ws.Cells.Item(1107, 6).Value = 1
$ws.Cells.Item(1108, 4).Value = 14
$ws.Cells.Item(1108, 5).Value = "[b'only']"
$ws.Cells.Item(1108, 6).Value = 1
$ws.Cells.Item(1109, 4).Value = 15
$ws.Cells.Item(1109, 5).Value = "[b'the']"
$ws.Cells.Item(1109, 6).Value = 1
$ws.Cells.Item(1110, 4).Value = 16
$ws.Cells.Item(1110, 5).Value = "[b'information']"
$ws.Cells.Item(1110, 6).Value = 1
$ws.Cells.Item(1111, 4).Value = 17
$ws.Cells.Item(1111, 5).Value = "[b'the']"
$ws.Cells.Item(1111, 6).Value = 1
$ws.Cells.Item(1112, 4).Value = 18
$ws.Cells.Item(1112, 5).Value = "[b'interviewer']"
$ws.Cells.Item(1112, 6).Value = 1
$ws.Cells.Item(1113, 4).Value = 19
$ws.Cells.Item(1113, 5).Value = "[b'needs']"
$ws.Cells.Item(1113, 6).Value = 1
$ws.Cells.Item(1094, 3).Interior.Color = 5296274

# --- Sentence 53: row 1115 starts with C=53 (green style s=4), rows 1115-1134 (D=0..19) ---
$ws.Cells.Item(1115, 3).Value = 53
$ws.Cells.Item(1115, 4).Value = 0
$ws.Cells.Item(1115, 5).Value = "[b'please']"
$ws.Cells.Item(1115, 6).Value = 1
$ws.Cells.Item(1116, 4).Value = 1
$ws.Cells.Item(1116, 5).Value = "[b'take']"
$ws.Cells.Item(1116, 6).Value = 1
$ws.Cells.Item(1117, 4).Value = 2
$ws.Cells.Item(1117, 5).Value = "[b'a']"
$ws.Cells.Item(1117, 6).Value = 1
$ws.Cells.Item(1118, 4).Value = 3
$ws.Cells.Item(1118, 5).Value = "[b'few']"
$ws.Cells.Item(1118, 6).Value = 1
$ws.Cells.Item(1119, 4).Value = 4
$ws.Cells.Item(1119, 5).Value = "[b'moments']"
$ws.Cells.Item(1119, 6).Value = 1
$ws.Cells.Item(1120, 4).Value = 5
$ws.Cells.Item(1120, 5).Value = "[b'to']"
$ws.Cells.Item(1120, 6).Value = 1
$ws.Cells.Item(1121, 4).Value = 6
$ws.Cells.Item(1121, 5).Value = "[b'learn']"
$ws.Cells.Item(1121, 6).Value = 1
$ws.Cells.Item(1122, 4).Value = 7
$ws.Cells.Item(1122, 5).Value = "[b'more']"
$ws.Cells.Item(1122, 6).Value = 1
$ws.Cells.Item(1123, 4).Value = 8
$ws.Cells.Item(1123, 5).Value = "[b'about']"
$ws.Cells.Item(1123, 6).Value = 1
$ws.Cells.Item(1124, 4).Value = 9
$ws.Cells.Item(1124, 5).Value = "[b'the']"
$ws.Cells.Item(1124, 6).Value = 1
$ws.Cells.Item(1125, 4).Value = 10
$ws.Cells.Item(1125, 5).Value = "[b'brochures']"
$ws.Cells.Item(1125, 6).Value = 1
$ws.Cells.Item(1126, 4).Value = 11
$ws.Cells.Item(1126, 5).Value = "[b'and']"
$ws.Cells.Item(1126, 6).Value = 1
$ws.Cells.Item(1127, 4).Value = 12
$ws.Cells.Item(1127, 5).Value = "[b'sample']"
$ws.Cells.Item(1127, 6).Value = 1
$ws.Cells.Item(1128, 4).Value = 13
$ws.Cells.Item(1128, 5).Value = "[b'packs']"
$ws.Cells.Item(1128, 6).Value = 1
$ws.Cells.Item(1129, 4).Value = 14
$ws.Cells.Item(1129, 5).Value = "[b'that']"
$ws.Cells.Item(1129, 6).Value = 1
$ws.Cells.Item(1130, 4).Value = 15
$ws.Cells.Item(1130, 5).Value = "[b'are']"
$ws.Cells.Item(1130, 6).Value = 1
$ws.Cells.Item(1131, 4).Value = 16
$ws.Cells.Item(1131, 5).Value = "[b'also']"
$ws.Cells.Item(1131, 6).Value = 1
$ws.Cells.Item(1132, 4).Value = 17
$ws.Cells.Item(1132, 5).Value = "[b'available']"
$ws.Cells.Item(1132, 6).Value = 1
$ws.Cells.Item(1133, 4).Value = 18
$ws.Cells.Item(1133, 5).Value = "[b'to']"
$ws.Cells.Item(1133, 6).Value = 1
$ws.Cells.Item(1134, 4).Value = 19
$ws.Cells.Item(1134, 5).Value = "[b'you']"
$ws.Cells.Item(1134, 6).Value = 1
$ws.Cells.Item(1115, 3).Interior.Color = 5296274

# --- Sentence 54: row 1136 starts with C=54 (no style), rows 1136-1155 (D=0..19), no words ---
$ws.Cells.Item(1136, 3).Value = 54
$ws.Cells.Item(1136, 4).Value = 0
$ws.Cells.Item(1137, 4).Value = 1
$ws.Cells.Item(1138, 4).Value = 2
$ws.Cells.Item(1139, 4).Value = 3
$ws.Cells.Item(1140, 4).Value = 4
$ws.Cells.Item(1141, 4).Value = 5
$ws.Cells.Item(1142, 4).Value = 6
$ws.Cells.Item(1143, 4).Value = 7
$ws.Cells.Item(1144, 4).Value = 8
$ws.Cells.Item(1145, 4).Value = 9
$ws.Cells.Item(1146, 4).Value = 10
$ws.Cells.Item(1147, 4).Value = 11
$ws.Cells.Item(1148, 4).Value = 12
$ws.Cells.Item(1149, 4).Value = 13
$ws.Cells.Item(1150, 4).Value = 14
$ws.Cells.Item(1151, 4).Value = 15
$ws.Cells.Item(1152, 4).Value = 16
$ws.Cells.Item(1153, 4).Value = 17
$ws.Cells.Item(1154, 4).Value = 18
$ws.Cells.Item(1155, 4).Value = 19

# --- Update view: selection + scroll position ---
$excel.ActiveWindow.ScrollRow = 1125
$ws.Range("D1136:D1155").Select()
